# Update latest output (run 166)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 295.0942507500001
$schedule.Range("F2").Value = 6.505605175264553
$schedule.Range("E3").Value = 428.2000905
$schedule.Range("F3").Value = 28.32011180555556

# --- Sheet "Detailed" ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B19").Value = -6.38688

$detailed.Range("B20").Value = -6.90384

$detailed.Range("B21").Value = -6.61075
$detailed.Range("C21").Value = "historical"

$detailed.Range("B22").Value = -5.50985
$detailed.Range("C22").Value = "historical"

$detailed.Range("B23").Value = -6.38177

$detailed.Range("B24").Value = -8.67718

$detailed.Range("B25").Value = -7.56508

$detailed.Range("B26").Value = -10

$detailed.Range("B27").Value = -12.01

$detailed.Range("B28").Value = -14

$detailed.Range("B29").Value = -8.64123

$detailed.Range("B30").Value = -12.01

$detailed.Range("B31").Value = -14

$detailed.Range("B32").Value = -10.81131

$detailed.Range("B33").Value = -12.01

$detailed.Range("B34").Value = -22.67512

$detailed.Range("B35").Value = -14

$detailed.Range("B36").Value = -8.92727

$detailed.Range("B37").Value = -7.69708

$detailed.Range("B38").Value = 1.6704

$detailed.Range("B39").Value = 10.2051

$detailed.Range("B40").Value = 27.5939

$detailed.Range("B41").Value = 57.31

$detailed.Range("B43").Value = 55.33037

$detailed.Range("B45").Value = 57.01318

$detailed.Range("B46").Value = 43.6164

$detailed.Range("B49").Value = 56.91608
